$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/ub04-admit-bill-type"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsInclude = $wb.Worksheets.Item("Include from UB04 Admit Bill ")
$wsInclude.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/ub04-admit-bill-type"
